$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Starfish diagram")

$text22 = '1.  Creare o mantenere un ambiente di lavoro eccessivamente pressante che possa influenzare negativamente la serenità e la produttività del team.     
2.eccessivo controllo per ogni dettaglio, a scapito della fiducia nel lavoro svolto dagli altri membri del gruppo.   '
$text23 = '
1. Prendere decisioni unilaterali sul lavoro da assegnare, senza coinvolgere attivamente tutti i membri del team, limitando l’autonomia e la partecipazione                                                             2. Intervenire sul lavoro altrui senza confronto preventivo, generando confusione e rendendo vani gli sforzi già fatti.                                                      3.Insistere nel voler realizzare funzionalità non previste o non chiaramente definite'

$ws.Range("D5").Value = $text22
$ws.Range("F5").Value = $text23

$ws.Rows.Item(5).RowHeight = 395.4

$ws.Range("K5").Select()
